$wb = $excel.ActiveWorkbook

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5895.7314
$ws.Range("I32").Value = 3538.4312
$ws.Range("K32").Value = 3538.4312
$ws.Range("M32").Value = -3251.4312

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2249.3977
$ws.Range("I63").Value = 2257.5
$ws.Range("K63").Value = 2257.5
$ws.Range("M63").Value = -1571.5

# ARM row 64
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 19650
$ws.Range("J64").Value = 19650
$ws.Range("L64").Value = 19650
$ws.Range("N64").Value = -20146

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2249.3977
$ws.Range("I66").Value = 2257.5
$ws.Range("K66").Value = 11287.5
$ws.Range("M66").Value = -7855.5

# ARM row 67
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 19650
$ws.Range("J67").Value = 19650
$ws.Range("L67").Value = 19650
$ws.Range("N67").Value = -21366

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1451.8
$ws.Range("I102").Value = 1413.2222
$ws.Range("J102").Value = 1799
$ws.Range("K102").Value = 1413.2222
$ws.Range("L102").Value = 1799
$ws.Range("M102").Value = 208.7778000000001
$ws.Range("N102").Value = -5043

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1835.5
$ws.Range("I122").Value = 1575.8462
$ws.Range("J122").Value = 2398.0833
$ws.Range("K122").Value = 4727.5386
$ws.Range("L122").Value = 7194.249899999999
$ws.Range("M122").Value = -2277.5386
$ws.Range("N122").Value = -12094.2499

# CRP row 9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 110071
$ws.Range("J9").Value = 110071
$ws.Range("L9").Value = 110071
$ws.Range("N9").Value = -110407

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1935.5555
$ws.Range("I122").Value = 1011
$ws.Range("J122").Value = 2199.7144
$ws.Range("K122").Value = 3033
$ws.Range("L122").Value = 6599.1432
$ws.Range("M122").Value = -583
$ws.Range("N122").Value = -11499.1432

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 173.16667
$ws.Range("I38").Value = 62.88889
$ws.Range("J38").Value = 504
$ws.Range("K38").Value = 188.66667
$ws.Range("L38").Value = 1512
$ws.Range("M38").Value = 158.33333
$ws.Range("N38").Value = -2206

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3202.879
$ws.Range("I39").Value = 1848.75
$ws.Range("J39").Value = 3389.6553
$ws.Range("K39").Value = 5546.25
$ws.Range("L39").Value = 10168.9659
$ws.Range("M39").Value = -5252.25
$ws.Range("N39").Value = -10756.9659

# CUL row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 293
$ws.Range("I40").Value = 121.4
$ws.Range("J40").Value = 1151
$ws.Range("K40").Value = 485.6
$ws.Range("L40").Value = 4604
$ws.Range("M40").Value = -416.6
$ws.Range("N40").Value = -4742

# CUL row 41
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M41").ClearContents()
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15676

# CUL row 42
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 1340
$ws.Range("J42").Value = 1340
$ws.Range("L42").Value = 4020
$ws.Range("N42").Value = -5088

# CUL row 44
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 764.3333
$ws.Range("I44").Value = 146.5
$ws.Range("J44").Value = 2000
$ws.Range("K44").Value = 439.5
$ws.Range("L44").Value = 6000
$ws.Range("M44").Value = -41.5
$ws.Range("N44").Value = -6796

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3500.75
$ws.Range("I46").Value = 4003
$ws.Range("J46").Value = 3333.3333
$ws.Range("K46").Value = 12009
$ws.Range("L46").Value = 9999.999899999999
$ws.Range("M46").Value = -11918
$ws.Range("N46").Value = -10181.9999

# CUL row 47
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1050.25
$ws.Range("I47").Value = 400.6
$ws.Range("J47").Value = 1514.2858
$ws.Range("K47").Value = 1201.8
$ws.Range("L47").Value = 4542.857400000001
$ws.Range("M47").Value = -770.8000000000002
$ws.Range("N47").Value = -5404.857400000001

# CUL row 48
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 3800
$ws.Range("J48").Value = 3800
$ws.Range("L48").Value = 11400
$ws.Range("N48").Value = -11900

# CUL row 49
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 3725
$ws.Range("I49").Value = 3000
$ws.Range("J49").Value = 3966.6667
$ws.Range("K49").Value = 9000
$ws.Range("L49").Value = 11900.0001
$ws.Range("M49").Value = -8844
$ws.Range("N49").Value = -12212.0001

# CUL row 50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 324
$ws.Range("I50").Value = 103
$ws.Range("J50").Value = 766
$ws.Range("K50").Value = 309
$ws.Range("L50").Value = 2298
$ws.Range("M50").Value = 172
$ws.Range("N50").Value = -3260

# CUL row 51
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2555
$ws.Range("I51").Value = 1000
$ws.Range("J51").Value = 2999.2856
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 8997.856800000001
$ws.Range("M51").Value = -2540
$ws.Range("N51").Value = -9917.856800000001

# CUL row 53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 324
$ws.Range("I53").Value = 103
$ws.Range("J53").Value = 766
$ws.Range("K53").Value = 309
$ws.Range("L53").Value = 2298
$ws.Range("M53").Value = 172
$ws.Range("N53").Value = -3260

# CUL row 54
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3119.2856
$ws.Range("I54").Value = 2871.25
$ws.Range("J54").Value = 3450
$ws.Range("K54").Value = 8613.75
$ws.Range("L54").Value = 10350
$ws.Range("M54").Value = -8054.75
$ws.Range("N54").Value = -11468

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 28842
$ws.Range("I55").Value = 87093.664
$ws.Range("J55").Value = 2628.75
$ws.Range("K55").Value = 261280.992
$ws.Range("L55").Value = 7886.25
$ws.Range("M55").Value = -261103.992
$ws.Range("N55").Value = -8240.25

# CUL row 61
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M61").ClearContents()
$ws.Range("H61").Value = 586.3889
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 586.3889
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 1759.1667
$ws.Range("N61").Value = -2189.1667

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1566.762
$ws.Range("I126").Value = 1315.8334
$ws.Range("J126").Value = 1901.3334
$ws.Range("K126").Value = 3947.5002
$ws.Range("L126").Value = 5704.0002
$ws.Range("M126").Value = -1477.5002
$ws.Range("N126").Value = -10644.0002

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6814.9165
$ws.Range("I132").Value = 8565.375
$ws.Range("J132").Value = 3314
$ws.Range("K132").Value = 25696.125
$ws.Range("L132").Value = 9942
$ws.Range("M132").Value = -23166.125
$ws.Range("N132").Value = -15002

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1116.4595
$ws.Range("I126").Value = 705.9048
$ws.Range("J126").Value = 1655.3125
$ws.Range("K126").Value = 2117.7144
$ws.Range("L126").Value = 4965.9375
$ws.Range("M126").Value = 352.2856000000002
$ws.Range("N126").Value = -9905.9375
